# Applies the cryptos list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the literal text into the cell without Excel re-parsing
    # number-looking strings (e.g. "575.90") into floating point
    # numbers, and without leaving the cells style index changed.
    $range = $ws.Range($cellRef)
    $originalStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $originalStyle
}

Set-TextValue 'D2' '65.149.99'
$ws.Range('E2').Value = '  +1.64%  '
Set-TextValue 'D3' '3.180.08'
$ws.Range('E3').Value = '  +3.88%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue 'D5' '575.90'
$ws.Range('E5').Value = '  +2.71%  '
Set-TextValue 'D6' '151.23'
$ws.Range('E6').Value = '  +5.02%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue 'D8' '3.176.81'
$ws.Range('E8').Value = '  +3.80%  '
$ws.Range('E9').Value = '  +3.29%  '
$ws.Range('E10').Value = '  +5.09%  '
Set-TextValue 'D11' '6.21'
$ws.Range('E11').Value = '  +2.21%  '
Set-TextValue 'D12' '0.508'
$ws.Range('E12').Value = '  +4.75%  '
$ws.Range('E13').Value = '  +19.62%  '
Set-TextValue 'D14' '38.18'
$ws.Range('E14').Value = '  +7.58%  '
Set-TextValue 'D15' '3.704.42'
$ws.Range('E15').Value = '  +4.01%  '
Set-TextValue 'D16' '65.235.90'
$ws.Range('E16').Value = '  +1.72%  '
Set-TextValue 'D17' '3.186.72'
$ws.Range('E17').Value = '  +3.84%  '
Set-TextValue 'D18' '7.21'
$ws.Range('E18').Value = '  +6.81%  '
$ws.Range('E19').Value = '  +1.16%  '
Set-TextValue 'D20' '514.60'
$ws.Range('E20').Value = '  +7.63%  '
Set-TextValue 'D21' '14.94'
$ws.Range('E21').Value = '  +6.95%  '
Set-TextValue 'D22' '0.735'
$ws.Range('E22').Value = '  +7.93%  '
Set-TextValue 'D23' '15.51'
$ws.Range('E23').Value = '  +8.63%  '
$ws.Range('E24').Value = '  +4.12%  '
Set-TextValue 'D25' '85.17'
$ws.Range('E25').Value = '  +3.61%  '
Set-TextValue 'D27' '9.08'
$ws.Range('E27').Value = '  +12.78%  '
$ws.Range('E28').Value = '  +4.24%  '
Set-TextValue 'D29' '2.20'
$ws.Range('E29').Value = '  +7.81%  '
Set-TextValue 'D30' '28.16'
$ws.Range('E30').Value = '  +6.87%  '
Set-TextValue 'D31' '2.78'
$ws.Range('E31').Value = '  +13.76%  '
$ws.Range('E32').Value = '  +7.46%  '
$ws.Range('E33').Value = '  +0.06%  '
Set-TextValue 'D34' '6.33'
$ws.Range('E34').Value = '  +10.68%  '
Set-TextValue 'D35' '6.72'
$ws.Range('E35').Value = '  +7.38%  '
Set-TextValue 'D36' '55.82'
$ws.Range('E36').Value = '  +1.55%  '
Set-TextValue 'D37' '0.0896'
$ws.Range('E37').Value = '  +10.32%  '
Set-TextValue 'D38' '480.73'
$ws.Range('E38').Value = '  +8.13%  '
Set-TextValue 'D39' '3.13'
$ws.Range('E39').Value = '  +10.07%  '
$ws.Range('E40').Value = '  +3.87%  '
Set-TextValue 'D41' '3.139.44'
$ws.Range('E41').Value = '  +4.93%  '
Set-TextValue 'D42' '8.66'
$ws.Range('E42').Value = '  +4.87%  '
$ws.Range('E43').Value = '  +4.41%  '
Set-TextValue 'D44' '0.289'
$ws.Range('E44').Value = '  +10.40%  '
Set-TextValue 'D45' '2.47'
$ws.Range('E45').Value = '  +13.93%  '
Set-TextValue 'D46' '29.45'
$ws.Range('E46').Value = '  +5.64%  '
Set-TextValue 'D47' '0.0₃0610'
$ws.Range('E47').Value = '  +17.55%  '
Set-TextValue 'D48' '0.999'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  +2.22%  '
$ws.Range('E50').Value = '  +10.76%  '
Set-TextValue 'D51' '122.16'
$ws.Range('E51').Value = '  +2.98%  '
